$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5493173599243164
$ws.Range("B1").Value = 1.107569336891174
$ws.Range("C1").Value = 5.272389888763428
$ws.Range("D1").Value = 4.129925727844238
$ws.Range("E1").Value = 0.8787401914596558
